$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = "NSE:DEVYANI"
$ws.Range("C2").Value = "NSE:3IINFOLTD"
$ws.Range("D2").ClearContents()
$ws.Range("E2").Value = "NSE:ADANIGREEN"
$ws.Range("F2").ClearContents()

# Row 3
$ws.Range("B3").Value = "NSE:GUFICBIO"
$ws.Range("C3").Value = "NSE:AUROPHARMA"
$ws.Range("D3").ClearContents()
$ws.Range("E3").Value = "NSE:ASTRAL"
$ws.Range("F3").ClearContents()

# Row 4
$ws.Range("B4").Value = "NSE:IDFNIFTYET"
$ws.Range("C4").Value = "NSE:CANFINHOME"
$ws.Range("E4").Value = "NSE:BHEL"
$ws.Range("F4").ClearContents()

# Row 5
$ws.Range("B5").Value = "NSE:JUBLFOOD"
$ws.Range("C5").Value = "NSE:DCI"
$ws.Range("E5").Value = "NSE:CESC"

# Row 6
$ws.Range("B6").Value = "NSE:LINDEINDIA"
$ws.Range("C6").ClearContents()
$ws.Range("E6").Value = "NSE:CONCOR"

# Row 7
$ws.Range("B7").Value = "NSE:MAXESTATES"
$ws.Range("C7").ClearContents()
$ws.Range("E7").Value = "NSE:DIVISLAB"

# Row 8
$ws.Range("B8").Value = "NSE:NYKAA"
$ws.Range("C8").ClearContents()
$ws.Range("E8").Value = "NSE:FEDERALBNK"

# Row 9
$ws.Range("B9").Value = "NSE:POLYMED"
$ws.Range("C9").ClearContents()
$ws.Range("E9").Value = "NSE:HDFCAMC"

# Row 10
$ws.Range("B10").ClearContents()
$ws.Range("C10").ClearContents()
$ws.Range("E10").Value = "NSE:IEX"

# Row 11
$ws.Range("B11").ClearContents()
$ws.Range("C11").ClearContents()
$ws.Range("E11").Value = "NSE:INDIGO"

# Row 12
$ws.Range("B12").ClearContents()
$ws.Range("C12").ClearContents()
$ws.Range("E12").Value = "NSE:LTTS"

# Row 13
$ws.Range("C13").ClearContents()
$ws.Range("E13").Value = "NSE:MFSL"

# Row 14
$ws.Range("C14").ClearContents()
$ws.Range("E14").Value = "NSE:PIIND"

# Remove rows 15 and 16 (no longer present in the updated table)
$ws.Rows("15:16").Delete()
